# "End of the project" - fill in the remainder of the "Journal de travail"
# sheet with the final days of work, update the total, and tidy up the
# now-unused trailing blank rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")

# The sheet currently has 15 unused placeholder rows (44-58) before the
# "Total" row. Only 12 of them are needed for the remaining journal
# entries, so remove 3 blank rows first - everything below (Total row,
# signature placeholders) shifts up to its final position automatically,
# and any formula ranges that are untouched stay intact.
$ws.Range("A56:A58").EntireRow.Delete()

# --- New journal entries -------------------------------------------------
# Column B values are entered in this exact order (44, 46, 45, 47, 48, ...)
# to match how the shared-string table ends up ordered.

$ws.Cells.Item(44, 1).Value = 45782
$ws.Cells.Item(44, 2).Value = "J'ai commencé la matinée par remettre en place mon infrastructure"
$ws.Cells.Item(44, 3).Value = 1
$ws.Range("A44:C44").RowHeight = 28.5

$ws.Cells.Item(46, 1).Value = 45782
$ws.Cells.Item(46, 2).Value = "Mon professeur est passé au cours de la matinée pour discuter avec moi de ce qu'il me fallait faire durant la remédiation"
$ws.Cells.Item(46, 3).Value = 0.5
$ws.Range("A46:C46").RowHeight = 28.5

$ws.Cells.Item(45, 1).Value = 45782
$ws.Cells.Item(45, 2).Value = "J'ai repris mes documents et ait relu ma documentation afin de me remttre dedans"
$ws.Cells.Item(45, 3).Value = 0.5
$ws.Range("A45:C45").RowHeight = 28.5

$ws.Cells.Item(47, 1).Value = 45782
$ws.Cells.Item(47, 2).Value = "J'ai termincé la matinée en mettant en forme mon code ainsi que mon infrastructure"
$ws.Cells.Item(47, 3).Value = 2
$ws.Range("A47:C47").RowHeight = 28.5

$ws.Cells.Item(48, 1).Value = 45782
$ws.Cells.Item(48, 2).Value = "J'ai passer un bon moment de l'après midi à résoudres les problèmes que j'avais avec la réservation"
$ws.Cells.Item(48, 3).Value = 2
$ws.Range("A48:C48").RowHeight = 28.5

$ws.Cells.Item(49, 1).Value = 45782
$ws.Cells.Item(49, 2).Value = "J'ai continuer en affichant les joueurs sur l'interface"
$ws.Cells.Item(49, 3).Value = 1

$ws.Cells.Item(50, 1).Value = 45782
$ws.Cells.Item(50, 2).Value = "J'ai terminé la journée en hébergeant mon site et en commencant la résolution de quelques problèmes du à l'hébergement"
$ws.Cells.Item(50, 3).Value = 1
$ws.Range("A50:C50").RowHeight = 28.5

$ws.Cells.Item(51, 1).Value = 45783
$ws.Cells.Item(51, 2).Value = "J'ai commencé la journée en résoluant mes problèmes du à l'hébergement "
$ws.Cells.Item(51, 3).Value = 0.75
$ws.Range("A51:C51").RowHeight = 28.5

$ws.Cells.Item(52, 1).Value = 45783
$ws.Cells.Item(52, 2).Value = "J'ai rencontré des problèmes d'accent et de majuscule qui m'ont pris un certain temps à résoudre à l'aide de mon professeur"
$ws.Cells.Item(52, 3).Value = 1
$ws.Range("A52:C52").RowHeight = 28.5

$ws.Cells.Item(53, 1).Value = 45783
$ws.Cells.Item(53, 2).Value = "J'ai effectuer une démonstration de mon application auprès de mon professeur"
$ws.Cells.Item(53, 3).Value = 0.25
$ws.Range("A53:C53").RowHeight = 28.5

$ws.Cells.Item(54, 1).Value = 45783
$ws.Cells.Item(54, 2).Value = "J'ai corriger mon rapport de projet et l'ai finalisé"
$ws.Cells.Item(54, 3).Value = 0.5

$ws.Cells.Item(55, 1).Value = 45783
$ws.Cells.Item(55, 2).Value = "J'ai mis à jour le journal de travail"
$ws.Cells.Item(55, 3).Value = 0.5

# --- Also bump up the time spent on the two last days before remediation -
$ws.Cells.Item(41, 3).Value = 1.5
$ws.Cells.Item(42, 3).Value = 2

# --- Total row now covers the extended range -----------------------------
$ws.Cells.Item(56, 3).Formula = "=SUM(C8:C55)"

# --- View tidy-up ----------------------------------------------------------
[void]$ws.Range("B48").Select()
